$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyle = $ws.Range("D6").Style

$ws.Range('D2').Value = '25.733.15'
$ws.Range('D3').Value = '1.627.67'
$ws.Range('E4').Value = '  -0.65%  '
$ws.Range('D5').Value = '''214.32'
$ws.Range('D5').Style = $defaultStyle
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  -0.60%  '
$ws.Range('E8').Value = '  -0.96%  '
$ws.Range('E9').Value = '  -0.87%  '
$ws.Range('D10').Value = '''19.49'
$ws.Range('D10').Style = $defaultStyle
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('D11').Value = '''0.0790'
$ws.Range('D11').Style = $defaultStyle
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '1.852.96'
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('D14').Value = '1.627.66'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '''0.552'
$ws.Range('D15').Style = $defaultStyle
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('E16').Value = '  -1.16%  '
$ws.Range('D17').Value = '''62.76'
$ws.Range('D17').Style = $defaultStyle
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('D18').Value = '25.729.24'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = '''4.43'
$ws.Range('D20').Style = $defaultStyle
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('D21').Value = '''191.20'
$ws.Range('D21').Style = $defaultStyle
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('D26').Value = '''142.29'
$ws.Range('D26').Style = $defaultStyle
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('E27').Value = '  +3.15%  '
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').Value = '''15.47'
$ws.Range('D29').Style = $defaultStyle
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  +1.44%  '
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('D33').Value = '''3.22'
$ws.Range('D33').Style = $defaultStyle
$ws.Range('E33').Value = '  -1.00%  '
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '1.138.62'
$ws.Range('E37').Value = '  +2.95%  '
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('D39').Value = '''0.542'
$ws.Range('D39').Style = $defaultStyle
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').Value = '''0.0155'
$ws.Range('D40').Style = $defaultStyle
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('E41').Value = '  -0.61%  '
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').Value = '''100.70'
$ws.Range('D43').Style = $defaultStyle
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''5.54'
$ws.Range('D44').Style = $defaultStyle
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').Value = '''0.802'
$ws.Range('D45').Style = $defaultStyle
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = '1.762.65'
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '''0.417'
$ws.Range('D49').Style = $defaultStyle
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = '''1.44'
$ws.Range('D50').Style = $defaultStyle
$ws.Range('E50').Value = '  +5.29%  '
$ws.Range('E51').Value = '  -0.54%  '
